$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume 1h (column E) updates for the refreshed symbol list.
# Cells are forced to Text format ("@") before assignment so that the numeric-looking
# strings (prices/percentages) are preserved verbatim instead of being parsed into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "317.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.85%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.46%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.185"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.30%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08011"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.63%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.490"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.53%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.533"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.25%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.942"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.27%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.971"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.64%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9420"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.25%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1323"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "11.58%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1938"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.57%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09053"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.11%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03354"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09537"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.60%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001397"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.39%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006497"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "10.67%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.392"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.79%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3519"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.25%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.542"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "24.37%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1315"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.47%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2420"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.49%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04379"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.57%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001228"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.52%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004271"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-8.40%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001327"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.42%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003984"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.01%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02376"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.38%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05165"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.64%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007732"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.65%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1401"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.69%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008575"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.18%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002105"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.00%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008963"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.12%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006468"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.34%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.02%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002863"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-13.39%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001687"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "68.96%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.02%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.02%"
